$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 236713.67
$ws.Range("J17").Value = 245443.81
$ws.Range("L17").Value = 736331.4299999999
$ws.Range("N17").Value = -736667.4299999999

$ws.Range("H40").Value = 1473.75
$ws.Range("J40").Value = 802
$ws.Range("L40").Value = 802
$ws.Range("N40").Value = -1152

$ws.Range("H107").Value = 55558948
$ws.Range("I107").Value = 3449.3333
$ws.Range("K107").Value = 3449.3333
$ws.Range("M107").Value = -1529.3333

$ws.Range("H113").Value = 6652.3335
$ws.Range("J113").Value = 6983.8
$ws.Range("L113").Value = 6983.8
$ws.Range("N113").Value = -13491.8

$ws.Range("H132").Value = 1322.6123
$ws.Range("I132").Value = 1300.5897
$ws.Range("J132").Value = 1408.5
$ws.Range("K132").Value = 3901.7691
$ws.Range("L132").Value = 4225.5
$ws.Range("M132").Value = -1371.7691
$ws.Range("N132").Value = -9285.5

$ws.Range("H135").Value = 2482.8333
$ws.Range("I135").Value = 2324.5
$ws.Range("K135").Value = 20920.5
$ws.Range("M135").Value = -18385.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1976.3334
$ws.Range("I21").Value = 1976.3334
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1976.3334
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -1602.3334
$ws.Range("N21").ClearContents()

$ws.Range("H32").Value = 224032.9
$ws.Range("I32").Value = 271027.62
$ws.Range("K32").Value = 271027.62
$ws.Range("M32").Value = -270740.62

$ws.Range("H61").Value = 2904266
$ws.Range("I61").Value = 71472.625
$ws.Range("J61").Value = 11969205
$ws.Range("K61").Value = 71472.625
$ws.Range("L61").Value = 11969205
$ws.Range("M61").Value = -71260.625
$ws.Range("N61").Value = -11969629

$ws.Range("H136").Value = 2904266
$ws.Range("I136").Value = 71472.625
$ws.Range("J136").Value = 11969205
$ws.Range("K136").Value = 214417.875
$ws.Range("L136").Value = 35907615
$ws.Range("M136").Value = -211867.875
$ws.Range("N136").Value = -35912715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6155.7856
$ws.Range("I86").Value = 3076.2222
$ws.Range("K86").Value = 3076.2222
$ws.Range("M86").Value = -1953.2222

$ws.Range("H89").Value = 6155.7856
$ws.Range("I89").Value = 3076.2222
$ws.Range("K89").Value = 15381.111
$ws.Range("M89").Value = -9765.111000000001

$ws.Range("H94").Value = 2387.2307
$ws.Range("I94").Value = 1460.25
$ws.Range("K94").Value = 1460.25
$ws.Range("M94").Value = -1009.25

$ws.Range("H126").Value = 91000
$ws.Range("J126").Value = 91000
$ws.Range("L126").Value = 91000
$ws.Range("N126").Value = -100880

$ws.Range("H134").Value = 29034346
$ws.Range("I134").Value = 1915.5238
$ws.Range("K134").Value = 5746.5714
$ws.Range("M134").Value = -3211.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2556.1042
$ws.Range("I31").Value = 3068.3809
$ws.Range("J31").Value = 2157.6667
$ws.Range("K31").Value = 3068.3809
$ws.Range("L31").Value = 2157.6667
$ws.Range("M31").Value = -2773.3809
$ws.Range("N31").Value = -2747.6667

$ws.Range("H34").Value = 2556.1042
$ws.Range("I34").Value = 3068.3809
$ws.Range("J34").Value = 2157.6667
$ws.Range("K34").Value = 3068.3809
$ws.Range("L34").Value = 2157.6667
$ws.Range("M34").Value = -2866.3809
$ws.Range("N34").Value = -2561.6667

$ws.Range("H60").Value = 45000
$ws.Range("J60").Value = 45000
$ws.Range("L60").Value = 45000
$ws.Range("N60").Value = -46022

$ws.Range("H74").Value = 42450
$ws.Range("J74").Value = 42450
$ws.Range("L74").Value = 42450
$ws.Range("N74").Value = -44198

$ws.Range("H77").Value = 42450
$ws.Range("J77").Value = 42450
$ws.Range("L77").Value = 127350
$ws.Range("N77").Value = -136086

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 4000
$ws.Range("I59").Value = 4000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 4000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -3417
$ws.Range("N59").ClearContents()

$ws.Range("H102").Value = 500000000
$ws.Range("I102").Value = 500000000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 500000000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -499998378
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 2671.75
$ws.Range("I126").Value = 2404
$ws.Range("K126").Value = 7212
$ws.Range("M126").Value = -4742

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H61").Value = 2331.3845
$ws.Range("I61").Value = 2331.6365
$ws.Range("K61").Value = 2331.6365
$ws.Range("M61").Value = -2129.6365

$ws.Range("H113").Value = 2331.3845
$ws.Range("I113").Value = 2331.6365
$ws.Range("K113").Value = 2331.6365
$ws.Range("M113").Value = -161.6365000000001

$ws.Range("H122").Value = 3144.6553
$ws.Range("J122").Value = 3939.2856
$ws.Range("L122").Value = 11817.8568
$ws.Range("N122").Value = -16717.8568

$ws.Range("H132").Value = 3623.182
$ws.Range("I132").Value = 3663.6
$ws.Range("J132").Value = 3536.5715
$ws.Range("K132").Value = 10990.8
$ws.Range("L132").Value = 10609.7145
$ws.Range("M132").Value = -8460.799999999999
$ws.Range("N132").Value = -15669.7145

$ws.Range("H134").Value = 69997.25
$ws.Range("J134").Value = 69997.25
$ws.Range("L134").Value = 69997.25
$ws.Range("N134").Value = -80137.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 15070
$ws.Range("I51").Value = 15070
$ws.Range("K51").Value = 15070
$ws.Range("M51").Value = -14560

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H107").Value = 2857801
$ws.Range("I107").Value = 826
$ws.Range("K107").Value = 2478
$ws.Range("M107").Value = -558

$ws.Range("H113").Value = 725.625
$ws.Range("I113").Value = 129
$ws.Range("J113").Value = 810.8570999999999
$ws.Range("K113").Value = 387
$ws.Range("L113").Value = 2432.5713
$ws.Range("M113").Value = 1783
$ws.Range("N113").Value = -6772.5713

$ws.Range("H123").Value = 84994.2
$ws.Range("J123").Value = 84994.2
$ws.Range("L123").Value = 84994.2
$ws.Range("N123").Value = -94794.2

$ws.Range("H132").Value = 2583.25
$ws.Range("I132").Value = 2263.0476
$ws.Range("K132").Value = 6789.1428
$ws.Range("M132").Value = -4259.1428
